# Commit: "Added chronic user chart"
# Updates the HeroinUseByFrequency sheet to show Occasional vs Chronic
# heroin users (National, Thousands) instead of the old "days in past
# month" frequency breakdown, adds a matching entry to the INFO tab
# index, and adjusts the view/selection state of the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. INFO sheet: add a row pointing at the new chart/sheet
# ---------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("INFO")
$wsInfo.Range("A4").Value = "HerionUseByFrequency"
$wsInfo.Range("B4").Value = "National, Thousands"

# ---------------------------------------------------------------
# 2. HeroinOverdoses sheet: no longer the active tab, cursor moved
# ---------------------------------------------------------------
$wsOver = $wb.Worksheets.Item("HeroinOverdoses")
$wsOver.Activate() | Out-Null
$wsOver.Range("D21").Select() | Out-Null

# ---------------------------------------------------------------
# 3. US_HeroinPast30Days sheet: cursor moved
# ---------------------------------------------------------------
$wsPast30 = $wb.Worksheets.Item("US_HeroinPast30Days")
$wsPast30.Activate() | Out-Null
$wsPast30.Range("B3").Select() | Out-Null

# ---------------------------------------------------------------
# 4. HeroinUseByFrequency sheet: replace the 4 frequency-bucket
#    columns (D,E - including the helper sum formula) with a
#    simpler 2-column Occasional/Chronic users breakdown.
# ---------------------------------------------------------------
$wsFreq = $wb.Worksheets.Item("HeroinUseByFrequency")
$wsFreq.Activate() | Out-Null

# Drop the old "21+ / 11-20 / 4-10 days" + Chronic-sum columns entirely.
$wsFreq.Range("D1:E12").Delete() | Out-Null

# New headers
$wsFreq.Range("B1").Value = "Occasional Users"
$wsFreq.Range("C1").Value = "Chronic Users (more than 10 days per/mo)"

# New data (National, Thousands)
$data = @(
    @(2000, 170, 1400),
    @(2001, 130, 1400),
    @(2002, 210, 1300),
    @(2003, 130, 1300),
    @(2004, 120, 1300),
    @(2005, 180, 1200),
    @(2006, 380, 1200),
    @(2007, 150, 1200),
    @(2008, 240, 1300),
    @(2009, 340, 1500),
    @(2010, 330, 1500)
)

$r = 2
foreach ($row in $data) {
    $wsFreq.Range("A$r").Value = $row[0]
    $wsFreq.Range("B$r").Value = $row[1]
    $wsFreq.Range("C$r").Value = $row[2]
    $r = $r + 1
}

# Widen column B for the new header text, grow the header row for wrapping
$wsFreq.Columns.Item(2).ColumnWidth = 13.498697916666666
$wsFreq.Rows.Item(1).RowHeight = 75

# View: zoomed to 125% with C1 selected
$win = $excel.ActiveWindow
$win.Zoom = 125
$wsFreq.Range("C1").Select() | Out-Null

# ---------------------------------------------------------------
# 5. Make INFO the active tab again (matches tabSelected moving
#    back to INFO, with the cursor left under the new row).
# ---------------------------------------------------------------
$wsInfo.Activate() | Out-Null
$wsInfo.Range("B5").Select() | Out-Null
